{"js": "// The only substantive content change in this revision is in the Belbin\n// roles table: Michaela's row (\"The Team Role Inventory Test\" and\n// \"The reality\" columns) changes from \"Everything balanced\" to\n// \"Plant, everything else balanced\". (Everything else in the original\n// diff is purely incidental run re-splitting / grammar-check proofErr\n// markers introduced by Word's editor and does not change the visible\n// text.)\nconst results = context.document.body.search(\"Everything balanced\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Plant, everything else balanced\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The only substantive content change in this revision is in the Belbin\n# roles table: Michaela's row (\"The Team Role Inventory Test\" and\n# \"The reality\" columns) changes from \"Everything balanced\" to\n# \"Plant, everything else balanced\". (Everything else in the original\n# diff is purely incidental run re-splitting / grammar-check proofErr\n# markers introduced by Word's editor and does not change the visible\n# text.)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Everything balanced\"\n$find.Replacement.Text = \"Plant, everything else balanced\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
